# Insert two new data rows at row 316 (pushing existing rows 316+ down by 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(316).Insert()
$ws.Rows.Item(316).Insert()

# New row 316: Hass / Primera
$ws.Cells.Item(316, 1).Value = 11
$ws.Cells.Item(316, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(316, 3).Value = "Bíobío"
$ws.Cells.Item(316, 4).Value = 44551
$ws.Cells.Item(316, 5).Value = 8
$ws.Cells.Item(316, 6).Value = "Fruta"
$ws.Cells.Item(316, 7).Value = 100106
$ws.Cells.Item(316, 8).Value = "Oleaginosos"
$ws.Cells.Item(316, 9).Value = 100106002
$ws.Cells.Item(316, 10).Value = "Palta"
$ws.Cells.Item(316, 11).Value = "Hass"
$ws.Cells.Item(316, 12).Value = "Primera"
$ws.Cells.Item(316, 13).Value = 270
$ws.Cells.Item(316, 14).Value = 3000
$ws.Cells.Item(316, 15).Value = 3500
$ws.Cells.Item(316, 16).Value = 3278
$ws.Cells.Item(316, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(316, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(316, 19).Value = 3278
$ws.Cells.Item(316, 20).Value = 1

# New row 317: Hass / Segunda
$ws.Cells.Item(317, 1).Value = 11
$ws.Cells.Item(317, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(317, 3).Value = "Bíobío"
$ws.Cells.Item(317, 4).Value = 44551
$ws.Cells.Item(317, 5).Value = 8
$ws.Cells.Item(317, 6).Value = "Fruta"
$ws.Cells.Item(317, 7).Value = 100106
$ws.Cells.Item(317, 8).Value = "Oleaginosos"
$ws.Cells.Item(317, 9).Value = 100106002
$ws.Cells.Item(317, 10).Value = "Palta"
$ws.Cells.Item(317, 11).Value = "Hass"
$ws.Cells.Item(317, 12).Value = "Segunda"
$ws.Cells.Item(317, 13).Value = 300
$ws.Cells.Item(317, 14).Value = 2500
$ws.Cells.Item(317, 15).Value = 2500
$ws.Cells.Item(317, 16).Value = 2500
$ws.Cells.Item(317, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(317, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(317, 19).Value = 2500
$ws.Cells.Item(317, 20).Value = 1
